$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the adjacent "sum" header (G1) to the new "Save" header (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" values for each data row (rows 2-25)
$saveValues = @(0,0,1,1,0,1,0,0,0,1,1,0,1,0,0,1,0,0,0,0,1,1,0,1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
